$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.146.31'
$ws.Range("E2").Value = '  -1.43%  '
$ws.Range("D3").Value = '1.838.01'
$ws.Range("E3").Value = '  -1.55%  '
$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.23%  '
$ws.Range("D5").Value = "'240.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.58%  '
$ws.Range("D6").Value = "'0.6845"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.92%  '
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("E8").Value = '  -3.28%  '
$ws.Range("D9").Value = "'0.07415"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -4.77%  '
$ws.Range("D10").Value = "'23.12"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.38%  '
$ws.Range("D11").Value = "'0.07645"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.49%  '
$ws.Range("D12").Value = '1.840.70'
$ws.Range("E12").Value = '  -1.20%  '
$ws.Range("D13").Value = "'5.056"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.55%  '
$ws.Range("D14").Value = "'0.6799"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.81%  '
$ws.Range("D15").Value = "'87.47"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.58%  '
$ws.Range("D16").Value = "'6.153"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -7.58%  '
$ws.Range("D17").Value = '29.125.65'
$ws.Range("E17").Value = '  -1.14%  '
$ws.Range("D18").Value = "'0.000008150"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.01%  '
$ws.Range("D19").Value = '2.077.64'
$ws.Range("E19").Value = '  -1.21%  '
$ws.Range("D20").Value = "'228.84"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -6.36%  '
$ws.Range("E21").Value = '  -2.39%  '
$ws.Range("D22").Value = "'0.9997"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.08%  '
$ws.Range("D23").Value = "'7.334"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.63%  '
$ws.Range("E24").Value = '  +0.22%  '
$ws.Range("D25").Value = "'159.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.08%  '
$ws.Range("E26").Value = '  -5.97%  '
$ws.Range("D27").Value = "'8.705"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.81%  '
$ws.Range("D28").Value = "'18.06"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.30%  '
$ws.Range("D29").Value = "'1.510"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.02%  '
$ws.Range("E30").Value = '  -0.20%  '
$ws.Range("D31").Value = "'4.137"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.27%  '
$ws.Range("E32").Value = '  -0.98%  '
$ws.Range("D33").Value = "'0.05261"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.83%  '
$ws.Range("D34").Value = "'0.7548"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -5.46%  '
$ws.Range("D35").Value = "'1.849"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -4.50%  '
$ws.Range("D36").Value = "'1.132"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.25%  '
$ws.Range("D37").Value = "'2.682"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.40%  '
$ws.Range("D38").Value = '1.291.97'
$ws.Range("E38").Value = '  -3.23%  '
$ws.Range("E39").Value = '  -3.06%  '
$ws.Range("D40").Value = "'2.718"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.47%  '
$ws.Range("D41").Value = "'0.9367"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.71%  '
$ws.Range("D42").Value = "'5.926"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.85%  '
$ws.Range("D43").Value = "'104.56"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.17%  '
$ws.Range("D44").Value = "'0.9996"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.16%  '
$ws.Range("E45").Value = '  +2.44%  '
$ws.Range("D46").Value = '1.979.18'
$ws.Range("E46").Value = '  -1.14%  '
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("D48").Value = "'64.58"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.02%  '
$ws.Range("D49").Value = "'9.482"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.69%  '
$ws.Range("E50").Value = '  -1.40%  '
$ws.Range("D51").Value = "'0.07466"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +17.66%  '
